$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns (in order) that receive new values for each data row: C,D,E,F,G,H,J,K,L,N,O
$cols = @("C","D","E","F","G","H","J","K","L","N","O")

# Each inner array: row number followed by the 11 new values (matching $cols order)
$data = @(
    @(2, 0.4242059179547937, 0.1934919952619367, 0.1847061095799241, 1.569126101909163, 0.9162775796953326, 0.9664544584577754, 0.2220857320067822, 1.389556250315991, 0.1675501044385541, 1.286538529941922, 3.803466808656538),
    @(3, 0.4193116607325891, 0.1903551570391926, 0.183556004613564, 1.572856189925844, 0.9192893309116528, 0.9723775012726463, 0.2219464923611625, 1.265718408431127, 0.167146121090763, 1.285202279938034, 3.821934823871459),
    @(4, 0.4164961663937135, 0.1884929830112156, 0.1829205418578894, 1.575963697964795, 0.9217119713380697, 0.9764365902572507, 0.2219527906590457, 1.189708508205314, 0.1669507612362295, 1.284826507135335, 3.835362702607085),
    @(5, 0.4153966376714635, 0.1877502867354508, 0.182679418098985, 1.577435602051182, 0.9228433157135356, 0.978196953668089, 0.2219784752563072, 1.158743051111486, 0.1668844412293637, 1.284785584153994, 3.841359747022466),
    @(6, 0.4152169531373602, 0.1876279411572597, 0.1826404582742995, 1.577692428650678, 0.9230398751469622, 0.9784956796898854, 0.2219841374993905, 1.153601886486172, 0.166874232806201, 1.284785577466693, 3.842387262579507),
    @(7, 0.4164811440159042, 0.1884829012097384, 0.1829172177051781, 1.575982716140018, 0.9217266457149833, 0.9764599008507844, 0.2219530433986066, 1.189290856816001, 0.1669498129518132, 1.28482550039368, 3.835441454949489),
    @(8, 0.4224791032751085, 0.192397216892644, 0.1842949092071535, 1.570242663116233, 0.917196963319256, 0.9684091052352954, 0.2220186918763716, 1.346852748435367, 0.1673999004436482, 1.285985722477903, 3.809401070009073),
    @(9, 0.4357407181945518, 0.2005761081359623, 0.1875555672110139, 1.565469025537709, 0.9128690794626237, 0.9559702319344723, 0.2228744060517656, 1.655955551951251, 0.1686988996751708, 1.291773857669739, 3.77491326624002),
    @(10, 0.4463933917397469, 0.2068874754244945, 0.1902897029258455, 1.565913810592605, 0.912474452132912, 0.9488703853088083, 0.2239448196654195, 1.883030467533445, 0.1699051154265021, 1.298149399008409, 3.759694034043719),
    @(11, 0.4514360961225918, 0.2098234740445832, 0.1916065917210652, 1.566974458462667, 0.912901558700554, 0.9460827549424096, 0.2245274084705429, 1.986308224748484, 0.1705081338457575, 1.301507100064342, 3.754970745686364),
    @(12, 0.4533738187837457, 0.2109445085899324, 0.1921157307508778, 1.567499497309498, 0.913150653322603, 0.9450906875177907, 0.2247617441997747, 2.02541185073602, 0.1707442516788618, 1.302844022827941, 3.753498712744147),
    @(13, 0.4529552452265477, 0.2107026645892915, 0.1920056140905828, 1.56738093327499, 0.913093119030151, 0.9453015211513218, 0.2247106660382414, 2.016990467211258, 0.1706930546553664, 1.302553187795809, 3.753801658269367),
    @(14, 0.4515949501225691, 0.2099155177127869, 0.1916482695241974, 1.567015180887196, 0.9129203003447088, 0.945999863535107, 0.2245464125898167, 1.989525426941327, 0.1705274039577986, 1.301615780053652, 3.754843294816141),
    @(15, 0.4507653934205109, 0.2094345672822442, 0.1914307465468355, 1.5668072157045, 0.9128258243312786, 0.946435893328669, 0.224447588648971, 1.972701514617427, 0.1704269485063534, 1.301050102109571, 3.755522560875761),
    @(16, 0.4460677869391816, 0.206696898876416, 0.190205107971785, 1.565861763669361, 0.9124587574055028, 0.9490614576814664, 0.2239086680891944, 1.876280405055809, 0.1698667959885611, 1.297939140519688, 3.760047000332321),
    @(17, 0.4432362556556768, 0.2050339853955023, 0.189471908142675, 1.565501572906868, 0.912389035278764, 0.9507853742663173, 0.2236025333442981, 1.817122279758905, 0.1695370400366301, 1.296147579289581, 3.763386227110942),
    @(18, 0.441626160694284, 0.2040836386970994, 0.1890570762098491, 1.565375190985776, 0.9124060237937357, 0.9518185421603675, 0.2234354579321192, 1.783094384377705, 0.1693524870143577, 1.295160196778184, 3.765513925258091),
    @(19, 0.4410841962585152, 0.203762921331105, 0.1889178053812799, 1.565346276952098, 0.9124215773464073, 0.9521755036656145, 0.223380436716738, 1.771572909278177, 0.1692908801005828, 1.294833294999364, 3.766269884069487),
    @(20, 0.4435357604840249, 0.2052103730028279, 0.1895492463012154, 1.565531554551001, 0.912390547708597, 0.9505975534465563, 0.2236341901437697, 1.823419958485147, 0.1695716141520549, 1.296333838338882, 3.763009328893588),
    @(21, 0.4519937385719857, 0.2101464719809627, 0.1917529467709187, 1.567119262577549, 0.9129686894873146, 0.9457930190152553, 0.2245942856584833, 1.997592742613392, 0.170575849093396, 1.301889346274777, 3.754528747487569),
    @(22, 0.4576855924248662, 0.2134262834588156, 0.1932541477136773, 1.568876145086108, 0.9138557897984469, 0.9430233692465322, 0.2253017307322409, 2.111391947171057, 0.1712774326174298, 1.305901437705131, 3.750831530649407),
    @(23, 0.4546327736923388, 0.2116708982063358, 0.192447368361151, 1.567872664796681, 0.9133356894681981, 0.944467702252652, 0.2249168477647032, 2.050659018985925, 0.1708988567421841, 1.303725337013859, 3.752635886238721),
    @(24, 0.4434002989285091, 0.2051306104522155, 0.1895142608638309, 1.565517748481781, 0.9123896861698739, 0.9506823362037409, 0.2236198503009561, 1.820572831002437, 0.1695559675298739, 1.296249497868089, 3.763179076892129),
    @(25, 0.4319930553073732, 0.1983101005878893, 0.1866138606205467, 1.566066413070146, 0.9135513994524416, 0.9589770313698835, 0.2225652364883715, 1.572331972786515, 0.168303116363397, 1.289833877522852, 3.782467073059848)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$colLetter$row").Value = $value
    }
}

$wb.Save()
